$wb = $excel.ActiveWorkbook

# Add the new "demand" sheet (the model now links company, hubs and demand)
# and place it after "hubs" (i.e. at the end of the tab strip).
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "demand"
$newSheet.Move($null, $wb.Worksheets.Item("hubs"))

# Re-acquire a fresh reference to the sheet by name before writing to it.
$ws = $wb.Worksheets.Item("demand")

# Fill in the "total demand" and "single hub capacity" rows first (label,
# then unit), then add the "goal" row.
$ws.Range("A1").Value = "total demand"
$ws.Range("A2").Value = "single hub capacity"
$ws.Range("C1").Value = "Box"
$ws.Range("C2").Value = "Box/Hub"
$ws.Range("A3").Value = "goal"
$ws.Range("C3").Value = "Dmnl"

# Column B: values
$ws.Range("B1").Value = 10000
$ws.Range("B2").Value = 1000
$ws.Range("B3").Value = 0.6

# Make "demand" the active sheet with A3 selected, as in the saved workbook.
$ws.Activate()
$ws.Range("A3").Select()
